$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (tab/workbook "(4)" -> "(5)")
$ws.Name = "repayment_20250915_20250915 (5)"

# --- Numeric cell updates ---
$ws.Range("D2").Value = 4
$ws.Range("H2").Value = 376
$ws.Range("H3").Value = 380
$ws.Range("D4").Value = 7
$ws.Range("H4").Value = 870
$ws.Range("J4").Value = 2
$ws.Range("H5").Value = 1.107
$ws.Range("H6").Value = 992
$ws.Range("H7").Value = 451
$ws.Range("H8").Value = 947
$ws.Range("D9").Value = 5
$ws.Range("H9").Value = 473
$ws.Range("D10").Value = 3
$ws.Range("H10").Value = 396
$ws.Range("H11").Value = 690
$ws.Range("H12").Value = 897
$ws.Range("H13").Value = 1.0569999999999999
$ws.Range("H14").Value = 412
$ws.Range("H15").Value = 465
$ws.Range("D16").Value = 8
$ws.Range("H16").Value = 1.1140000000000001
$ws.Range("J16").Value = 3
$ws.Range("H17").Value = 694
$ws.Range("D18").Value = 3
$ws.Range("H18").Value = 1.665

# --- Text cell updates (values that look numeric but must stay text) ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "554,687.00"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.36"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5,841,430.00"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "3.53"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "13.54"
$ws.Range("L4").NumberFormat = "@"
$ws.Range("L4").Value = "6.45"
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "1.51"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "484,839.00"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.26"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "529,441.00"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.31"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,014,658.00"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "1.56"
$ws.Range("K16").NumberFormat = "@"
$ws.Range("K16").Value = "7.37"
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "10.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "924,590.00"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "0.91"

# Strip the temporary text-format styling so cells end up with the default style (matches source)
$ws.Range("E2").ClearFormats()
$ws.Range("G2").ClearFormats()
$ws.Range("E4").ClearFormats()
$ws.Range("G4").ClearFormats()
$ws.Range("K4").ClearFormats()
$ws.Range("L4").ClearFormats()
$ws.Range("K8").ClearFormats()
$ws.Range("E9").ClearFormats()
$ws.Range("G9").ClearFormats()
$ws.Range("E10").ClearFormats()
$ws.Range("G10").ClearFormats()
$ws.Range("E16").ClearFormats()
$ws.Range("G16").ClearFormats()
$ws.Range("K16").ClearFormats()
$ws.Range("L16").ClearFormats()
$ws.Range("E18").ClearFormats()
$ws.Range("G18").ClearFormats()
